$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.494.92"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.107.32"
$ws.Range("E3").Value = "  +1.43%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.11"
$ws.Range("E5").Value = "  +1.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.22"
$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.105.51"
$ws.Range("E8").Value = "  +1.45%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.23"
$ws.Range("E10").Value = "  -0.88%  "

$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("E12").Value = "  +1.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.641.07"
$ws.Range("E13").Value = "  +1.47%  "

$ws.Range("E14").Value = "  +0.92%  "

$ws.Range("E15").Value = "  +0.54%  "

$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.576.94"
$ws.Range("E17").Value = "  -0.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.104.99"
$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("E19").Value = "  +0.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.79"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("E21").Value = "  -0.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "336.51"
$ws.Range("E22").Value = "  +1.41%  "

$ws.Range("E24").Value = "  +2.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.58"
$ws.Range("E25").Value = "  +1.30%  "

$ws.Range("E26").Value = "  -0.53%  "

$ws.Range("E27").Value = "  +0.19%  "

$ws.Range("E28").Value = "  +1.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.53"
$ws.Range("E29").Value = "  +1.99%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.87"
$ws.Range("E31").Value = "  +2.33%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.94"
$ws.Range("E32").Value = "  +1.13%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "157.37"
$ws.Range("E34").Value = "  +1.82%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.65"
$ws.Range("E35").Value = "  +2.95%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.11"
$ws.Range("E36").Value = "  +2.46%  "

$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.96"
$ws.Range("E37").Value = "  -0.75%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.28"
$ws.Range("E38").Value = "  +1.16%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0661"
$ws.Range("E39").Value = "  -1.65%  "

$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.147.47"
$ws.Range("E40").Value = "  +1.38%  "

$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.686"
$ws.Range("E41").Value = "  +4.78%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.94"
$ws.Range("E42").Value = "  +0.60%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.51"
$ws.Range("E43").Value = "  +10.82%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "36.83"
$ws.Range("E44").Value = "  +0.84%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.302.48"
$ws.Range("E46").Value = "  +2.13%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0260"
$ws.Range("E47").Value = "  +0.56%  "

$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.977"
$ws.Range("E48").Value = "  +4.41%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.78"
$ws.Range("E49").Value = "  +0.84%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.02"
$ws.Range("E50").Value = "  +2.10%  "

$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.732"
$ws.Range("E51").Value = "  -0.06%  "
